{"js": "// Update the worksheet date header and the 25 \"two-digit \u00f7 one-digit\"\n// division problems/answers laid out in the 5-column practice table.\n// Each old value is unique within its own table cell, so we scope the\n// search to that cell's range (rather than doing a document-wide\n// search/replace) to avoid any cross-cell collisions between an old\n// value in one cell and a newly written value in another cell.\n\nasync function replaceInRange(range, oldText, newText) {\n  const results = range.search(oldText, { matchCase: true, matchWildcards: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Could not find text to replace: \" + oldText);\n  }\n\n  // Replace the first (and expected only) match, preserving the run's\n  // existing character formatting (font/size) since we are only\n  // swapping the text inside the matched range.\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 1) Date header paragraph above the table.\nawait replaceInRange(context.document.body.getRange(), \"2025-06-05 Thursday\", \"2025-06-06 Friday\");\n\n// 2) The practice table: 20 rows x 5 columns, with problems only in\n//    rows 0, 4, 8, 12, 16 (the other rows are spacer rows).\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nconst cellUpdates = [\n  // row, col, oldText,        newText\n  [0, 0, \"66\u00f77=9, 3\", \"53\u00f77=7, 4\"],\n  [0, 1, \"29\u00f73=9, 2\", \"56\u00f74=14, 0\"],\n  [0, 2, \"14\u00f78=1, 6\", \"19\u00f77=2, 5\"],\n  [0, 3, \"52\u00f72=26, 0\", \"86\u00f77=12, 2\"],\n  [0, 4, \"67\u00f79=7, 4\", \"91\u00f77=13, 0\"],\n\n  [4, 0, \"10\u00f79=1, 1\", \"92\u00f77=13, 1\"],\n  [4, 1, \"38\u00f75=7, 3\", \"85\u00f72=42, 1\"],\n  [4, 2, \"10\u00f72=5, 0\", \"53\u00f74=13, 1\"],\n  [4, 3, \"51\u00f77=7, 2\", \"96\u00f78=12, 0\"],\n  [4, 4, \"48\u00f73=16, 0\", \"23\u00f73=7, 2\"],\n\n  [8, 0, \"21\u00f78=2, 5\", \"20\u00f79=2, 2\"],\n  [8, 1, \"63\u00f79=7, 0\", \"96\u00f72=48, 0\"],\n  [8, 2, \"75\u00f76=12, 3\", \"58\u00f73=19, 1\"],\n  [8, 3, \"62\u00f74=15, 2\", \"94\u00f72=47, 0\"],\n  [8, 4, \"42\u00f73=14, 0\", \"19\u00f72=9, 1\"],\n\n  [12, 0, \"26\u00f73=8, 2\", \"85\u00f76=14, 1\"],\n  [12, 1, \"71\u00f76=11, 5\", \"29\u00f79=3, 2\"],\n  [12, 2, \"70\u00f74=17, 2\", \"15\u00f76=2, 3\"],\n  [12, 3, \"49\u00f74=12, 1\", \"67\u00f78=8, 3\"],\n  [12, 4, \"44\u00f75=8, 4\", \"30\u00f75=6, 0\"],\n\n  [16, 0, \"19\u00f77=2, 5\", \"23\u00f72=11, 1\"],\n  [16, 1, \"81\u00f76=13, 3\", \"77\u00f73=25, 2\"],\n  [16, 2, \"29\u00f78=3, 5\", \"51\u00f78=6, 3\"],\n  [16, 3, \"35\u00f75=7, 0\", \"94\u00f79=10, 4\"],\n  [16, 4, \"63\u00f77=9, 0\", \"41\u00f76=6, 5\"],\n];\n\nfor (const [row, col, oldText, newText] of cellUpdates) {\n  const cell = table.getCell(row, col);\n  const cellRange = cell.body.getRange();\n  await replaceInRange(cellRange, oldText, newText);\n}\n", "ps1": "# Update the worksheet date header and the 25 \"two-digit \u00f7 one-digit\"\n# division problems/answers in the 5-column practice table.\n#\n# Each cell's text is set directly via its own Range object\n# ($cell.Range.Text = ...) rather than Find/Replace. Several old\n# values re-appear as *new* values elsewhere in the table (e.g.\n# \"19\u00f77=2, 5\" is the new text for one cell and the old text of a\n# later cell), so a document-wide Find/Replace pass could re-match\n# text this script just wrote. Addressing each paragraph/cell range\n# directly avoids that ambiguity and keeps every other run's\n# formatting (font/size) untouched.\n\n$d = $word.ActiveDocument\n\n# 1) Date header paragraph above the table.\n$d.Paragraphs.Item(1).Range.Text = \"2025-06-06 Friday\"\n\n# 2) The practice table: 20 rows x 5 columns (1-based via COM), with\n#    problems only in rows 1, 5, 9, 13, 17 (the other rows are spacers).\n$tbl = $d.Tables.Item(1)\n\n$cellUpdates = @(\n    @(1, 1, \"53\u00f77=7, 4\"),\n    @(1, 2, \"56\u00f74=14, 0\"),\n    @(1, 3, \"19\u00f77=2, 5\"),\n    @(1, 4, \"86\u00f77=12, 2\"),\n    @(1, 5, \"91\u00f77=13, 0\"),\n\n    @(5, 1, \"92\u00f77=13, 1\"),\n    @(5, 2, \"85\u00f72=42, 1\"),\n    @(5, 3, \"53\u00f74=13, 1\"),\n    @(5, 4, \"96\u00f78=12, 0\"),\n    @(5, 5, \"23\u00f73=7, 2\"),\n\n    @(9, 1, \"20\u00f79=2, 2\"),\n    @(9, 2, \"96\u00f72=48, 0\"),\n    @(9, 3, \"58\u00f73=19, 1\"),\n    @(9, 4, \"94\u00f72=47, 0\"),\n    @(9, 5, \"19\u00f72=9, 1\"),\n\n    @(13, 1, \"85\u00f76=14, 1\"),\n    @(13, 2, \"29\u00f79=3, 2\"),\n    @(13, 3, \"15\u00f76=2, 3\"),\n    @(13, 4, \"67\u00f78=8, 3\"),\n    @(13, 5, \"30\u00f75=6, 0\"),\n\n    @(17, 1, \"23\u00f72=11, 1\"),\n    @(17, 2, \"77\u00f73=25, 2\"),\n    @(17, 3, \"51\u00f78=6, 3\"),\n    @(17, 4, \"94\u00f79=10, 4\"),\n    @(17, 5, \"41\u00f76=6, 5\")\n)\n\nforeach ($update in $cellUpdates) {\n    $row = $update[0]\n    $col = $update[1]\n    $newText = $update[2]\n    $tbl.Cell($row, $col).Range.Text = $newText\n}\n"}
